$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.971.04'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '1.639.77'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.53%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5083'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("E9").Value = '  -0.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.82'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07746'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.298'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("D13").Value = '1.644.65'
$ws.Range("E13").Value = '  -0.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5462'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.03%  '
$ws.Range("D15").Value = '0.0₅7734'
$ws.Range("E15").Value = '  -1.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.67%  '
$ws.Range("D17").Value = '25.988.92'
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.002'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.461'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '196.00'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.947'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.139'
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = '  -0.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.899'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.72'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1271'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +11.22%  '
$ws.Range("E27").Value = '  -0.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.63'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  -2.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.260'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.207'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.549'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.375'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9172'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.567'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.06%  '
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.133.39'
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5535'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.41%  '
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("E40").Value = '  -0.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.586'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.16%  '
$ws.Range("E42").Value = '  -1.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '98.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.50%  '
$ws.Range("D44").Value = '0.0₈121'
$ws.Range("E44").Value = '  -8.35%  '
$ws.Range("D45").Value = '1.774.40'
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4515'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.23'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9993'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05182'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.503'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.003'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.53%  '
